# Insert a new weekly price-report row at row 101 (pushes existing rows
# 101..178 down to 102..179) and fill it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 101..178 down by inserting a new row at position 101.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Cells.Item(101, 1).Value  = 11
$ws.Cells.Item(101, 2).Value  = 'Vega Monumental Concepción'
$ws.Cells.Item(101, 3).Value  = 'Bíobío'
$ws.Cells.Item(101, 4).Value  = 44586
$ws.Cells.Item(101, 5).Value  = 8
$ws.Cells.Item(101, 6).Value  = 100114013
$ws.Cells.Item(101, 7).Value  = 'Zanahoria'
$ws.Cells.Item(101, 8).Value  = 'Sin especificar'
$ws.Cells.Item(101, 9).Value  = 'Primera'
$ws.Cells.Item(101, 10).Value = 220
$ws.Cells.Item(101, 11).Value = 9000
$ws.Cells.Item(101, 12).Value = 9500
$ws.Cells.Item(101, 13).Value = 9273
$ws.Cells.Item(101, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(101, 15).Value = 'Chillán'
$ws.Cells.Item(101, 16).Value = 464
$ws.Cells.Item(101, 17).Value = 20
$ws.Cells.Item(101, 18).Value = 'Hortaliza'
